$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "-"
$ws.Range("F2").Value = "MCT-2A-Sistemas Digitais"

$ws.Range("C3").Value = "-"
$ws.Range("F3").Value = "MCT-2A-Sistemas Digitais"

$ws.Range("C4").Value = "MCT-3A-Laboratório de Eletroeletrônica"
$ws.Range("D4").Value = "-"
$ws.Range("F4").Value = "MCT-2A-Sistemas Digitais"

$ws.Range("C6").Value = "MCT-3A-Laboratório de Eletroeletrônica"
$ws.Range("D6").Value = "MCT-3A-Laboratório de Eletroeletrônica"
$ws.Range("F6").Value = "MCT-2A-Sistemas Digitais"

$ws.Range("C7").Value = "MCT-3A-Laboratório de Eletroeletrônica"
$ws.Range("D7").Value = "MCT-3A-Laboratório de Eletroeletrônica"
$ws.Range("F7").Value = "MCT-2A-Sistemas Digitais"

$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "MCT-3A-Laboratório de Eletroeletrônica"
$ws.Range("F8").Value = "MCT-2A-Sistemas Digitais"
